$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string must be forced to
# Text so Excel does not coerce them into numeric cells (which would drop
# formatting such as trailing zeros). NumberFormat is reset with ClearFormats
# afterwards so the cell keeps its original (default) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '67.597.38'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').Value = '2.516.58'
$ws.Range('E3').Value = '  +0.37%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '591.49'
$ws.Range('E5').Value = '  +2.16%  '
Set-TextValue $ws.Range('D6') '176.50'
$ws.Range('E6').Value = '  +5.73%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue $ws.Range('D8') '0.531'
$ws.Range('E8').Value = '  +2.26%  '
$ws.Range('D9').Value = '2.515.75'
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('E11').Value = '  +3.06%  '
Set-TextValue $ws.Range('D12') '5.16'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '2.976.67'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '67.416.55'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').Value = '2.520.31'
$ws.Range('E18').Value = '  +0.29%  '
Set-TextValue $ws.Range('D19') '7.95'
$ws.Range('E19').Value = '  +4.63%  '
Set-TextValue $ws.Range('D20') '11.40'
$ws.Range('E20').Value = '  +1.59%  '
Set-TextValue $ws.Range('D21') '360.87'
$ws.Range('E21').Value = '  +5.20%  '
$ws.Range('E22').Value = '  -0.34%  '
Set-TextValue $ws.Range('D23') '4.64'
$ws.Range('E23').Value = '  +2.22%  '
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('E25').Value = '  -0.06%  '
Set-TextValue $ws.Range('D26') '71.01'
$ws.Range('E26').Value = '  +3.18%  '
Set-TextValue $ws.Range('D27') '10.22'
$ws.Range('E27').Value = '  +3.24%  '
Set-TextValue $ws.Range('D28') '0.998'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').Value = '2.643.61'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '0.0₃0985'
$ws.Range('E30').Value = '  +0.84%  '
Set-TextValue $ws.Range('D31') '543.61'
$ws.Range('E31').Value = '  +4.04%  '
Set-TextValue $ws.Range('D32') '8.22'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('E34').Value = '  +3.00%  '
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +1.08%  '
Set-TextValue $ws.Range('D38') '155.67'
$ws.Range('E38').Value = '  -0.91%  '
Set-TextValue $ws.Range('D39') '18.78'
$ws.Range('E39').Value = '  +1.76%  '
Set-TextValue $ws.Range('D40') '18.62'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D44') '1.00'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D45') '2.52'
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('E46').Value = '  +1.17%  '
Set-TextValue $ws.Range('D47') '146.30'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('E51').Value = '  +0.21%  '
